$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 83
$ws.Range("I8").Value = 82.888885
$ws.Range("K8").Value = 248.666655
$ws.Range("M8").Value = -109.666655
$ws.Range("H9").Value = 868.53845
$ws.Range("I9").Value = 349.83334
$ws.Range("J9").Value = 1313.1428
$ws.Range("K9").Value = 349.83334
$ws.Range("L9").Value = 1313.1428
$ws.Range("M9").Value = -180.83334
$ws.Range("N9").Value = -1651.1428
$ws.Range("H15").Value = 643.3684
$ws.Range("I15").Value = 643.3684
$ws.Range("K15").Value = 1930.1052
$ws.Range("M15").Value = -1761.1052
$ws.Range("H33").Value = 866.6667
$ws.Range("I33").Value = 1080
$ws.Range("K33").Value = 1080
$ws.Range("M33").Value = -851
$ws.Range("H39").Value = 57.9
$ws.Range("I39").Value = 46.555557
$ws.Range("K39").Value = 139.666671
$ws.Range("M39").Value = 156.333329
$ws.Range("H64").Value = 5495.3335
$ws.Range("J64").Value = 5496.5
$ws.Range("L64").Value = 5496.5
$ws.Range("N64").Value = -5992.5
$ws.Range("H67").Value = 5495.3335
$ws.Range("J67").Value = 5496.5
$ws.Range("L67").Value = 5496.5
$ws.Range("N67").Value = -7212.5
$ws.Range("H100").Value = 1028
$ws.Range("I100").Value = 851.36365
$ws.Range("K100").Value = 851.36365
$ws.Range("M100").Value = -310.36365
$ws.Range("H137").Value = 3674.111
$ws.Range("I137").Value = 2331.3333
$ws.Range("K137").Value = 6993.999899999999
$ws.Range("M137").Value = -4443.999899999999
$ws.Range("H138").Value = 5802.2764
$ws.Range("I138").Value = 5714.636
$ws.Range("J138").Value = 5829.0557
$ws.Range("K138").Value = 17143.908
$ws.Range("L138").Value = 17487.1671
$ws.Range("M138").Value = -12003.908
$ws.Range("N138").Value = -27767.1671
$ws.Range("H141").Value = 2673.2354
$ws.Range("I141").Value = 2753.4375
$ws.Range("K141").Value = 8260.3125
$ws.Range("M141").Value = -3080.3125

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6971
$ws.Range("I28").Value = 6971
$ws.Range("K28").Value = 6971
$ws.Range("M28").Value = -6779
$ws.Range("H63").Value = 1999
$ws.Range("I63").Value = 1999
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1999
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1313
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1999
$ws.Range("I66").Value = 1999
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9995
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6563
$ws.Range("N66").ClearContents()
$ws.Range("H99").Value = 6971
$ws.Range("I99").Value = 6971
$ws.Range("K99").Value = 6971
$ws.Range("M99").Value = -3976
$ws.Range("H102").Value = 2795.3157
$ws.Range("I102").Value = 2857.0625
$ws.Range("K102").Value = 2857.0625
$ws.Range("M102").Value = -1235.0625
$ws.Range("H122").Value = 4996.625
$ws.Range("I122").Value = 4996.1665
$ws.Range("K122").Value = 14988.4995
$ws.Range("M122").Value = -12538.4995
$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5117.8
$ws.Range("I86").Value = 3113.5
$ws.Range("K86").Value = 3113.5
$ws.Range("M86").Value = -1990.5
$ws.Range("H89").Value = 5117.8
$ws.Range("I89").Value = 3113.5
$ws.Range("K89").Value = 15567.5
$ws.Range("M89").Value = -9951.5
$ws.Range("H94").Value = 847.5714
$ws.Range("I94").Value = 822.1667
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 822.1667
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -371.1667
$ws.Range("N94").Value = -1902
$ws.Range("H105").Value = 1702.4348
$ws.Range("I105").Value = 1327.5555
$ws.Range("K105").Value = 1327.5555
$ws.Range("M105").Value = 419.4445000000001
$ws.Range("H107").Value = 1968.8
$ws.Range("I107").Value = 1709.375
$ws.Range("K107").Value = 1709.375
$ws.Range("M107").Value = 210.625

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H86").Value = 7374.75
$ws.Range("J86").Value = 3500
$ws.Range("L86").Value = 3500
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 7374.75
$ws.Range("J89").Value = 3500
$ws.Range("L89").Value = 17500
$ws.Range("N89").Value = -28732
$ws.Range("H141").Value = 377626.62
$ws.Range("J141").Value = 377626.62
$ws.Range("L141").Value = 377626.62
$ws.Range("N141").Value = -387986.62

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1225.8334
$ws.Range("I117").Value = 1169
$ws.Range("J117").Value = 1305.4
$ws.Range("K117").Value = 3507
$ws.Range("L117").Value = 3916.2
$ws.Range("M117").Value = -65
$ws.Range("N117").Value = -10800.2
$ws.Range("H121").Value = 3141.1428
$ws.Range("I121").Value = 586.73334
$ws.Range("K121").Value = 1760.20002
$ws.Range("M121").Value = -450.20002

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1130.3334
$ws.Range("I22").Value = 1180
$ws.Range("J22").Value = 1080.6666
$ws.Range("K22").Value = 1180
$ws.Range("L22").Value = 1080.6666
$ws.Range("M22").Value = -651
$ws.Range("N22").Value = -2138.6666
$ws.Range("H70").Value = 5674.769
$ws.Range("I70").Value = 4499.5
$ws.Range("J70").Value = 5888.4546
$ws.Range("K70").Value = 4499.5
$ws.Range("L70").Value = 5888.4546
$ws.Range("M70").Value = -4229.5
$ws.Range("N70").Value = -6428.4546
$ws.Range("H73").Value = 5674.769
$ws.Range("I73").Value = 4499.5
$ws.Range("J73").Value = 5888.4546
$ws.Range("K73").Value = 4499.5
$ws.Range("L73").Value = 5888.4546
$ws.Range("M73").Value = -3563.5
$ws.Range("N73").Value = -7760.4546
$ws.Range("H96").Value = 19988.5
$ws.Range("J96").Value = 19988.5
$ws.Range("L96").Value = 19988.5
$ws.Range("N96").Value = -25480.5
$ws.Range("H97").Value = 990.5
$ws.Range("I97").Value = 972
$ws.Range("K97").Value = 972
$ws.Range("M97").Value = -476
$ws.Range("H113").Value = 1925
$ws.Range("I113").Value = 1925
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1925
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 245
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2680.9167
$ws.Range("I132").Value = 2561
$ws.Range("K132").Value = 7683
$ws.Range("M132").Value = -5153

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8256.429
$ws.Range("I136").Value = 10287.286
$ws.Range("J136").Value = 6225.5713
$ws.Range("K136").Value = 30861.858
$ws.Range("L136").Value = 18676.7139
$ws.Range("M136").Value = -28311.858
$ws.Range("N136").Value = -23776.7139

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1500
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4246
$ws.Range("H100").Value = 8334770.5
$ws.Range("I100").Value = 12500981
$ws.Range("K100").Value = 25001962
$ws.Range("M100").Value = -25001421
$ws.Range("H109").Value = 58376
$ws.Range("J109").Value = 58376
$ws.Range("L109").Value = 58376
$ws.Range("N109").Value = -61150
$ws.Range("H132").Value = 6736.8237
$ws.Range("I132").Value = 6501.7334
$ws.Range("K132").Value = 19505.2002
$ws.Range("M132").Value = -16975.2002
